# This script updates the 56 changed result cells in Sheet1 of the
# AV-MCPS process analysis workbook, reflecting the recalculated values
# from re-running the first three simulations ("analisis de las 3
# primeras simulaciones").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 48.11829044092123;
    "B3" = 63.89447360117288;
    "B4" = 72.16482847090539;
    "I5" = 96.1810605035658;
    "I6" = 96.21062246260344;
    "I7" = 96.20833922475407;
    "C8" = 95.74158889681613;
    "C9" = 92.24143140201321;
    "C10" = 94.79496498136005;
    "D11" = 99.03061139087326;
    "D12" = 98.54381242096044;
    "D13" = 99.41938885197247;
    "E14" = 98.82590060072793;
    "E15" = 98.95326730288227;
    "E16" = 98.89907547862148;
    "F17" = 98.43840310256745;
    "F18" = 98.50030733847247;
    "F19" = 98.43579522570374;
    "G20" = 97.97738793102236;
    "G21" = 98.05931384705194;
    "G22" = 98.02667089791233;
    "H23" = 97.57296703696768;
    "H24" = 97.58151027190115;
    "H25" = 97.56981416344483;
    "B26" = 75.32635977857164;
    "B27" = 76.39588681364611;
    "I28" = 96.09130061619031;
    "I29" = 96.1397816262419;
    "C30" = 94.73842084538811;
    "C31" = 94.4293374377917;
    "D32" = 98.89653647142185;
    "D33" = 99.1102226565378;
    "E34" = 98.97051661478353;
    "E35" = 98.93800044440172;
    "F36" = 98.48487032153372;
    "F37" = 98.55562060065448;
    "G38" = 97.97942797929666;
    "G39" = 97.99031180030244;
    "H40" = 97.58371768390084;
    "H41" = 97.56967667654615;
    "B42" = 70.6028786428225;
    "B43" = 88.74795946687171;
    "I44" = 96.14242532291524;
    "I45" = 96.21197634869935;
    "C46" = 95.11609911347131;
    "C47" = 93.26241810059329;
    "D48" = 98.90991227767221;
    "D49" = 99.05273299697252;
    "E50" = 98.87499492763729;
    "E51" = 98.34830869779729;
    "F52" = 98.5437223903388;
    "F53" = 98.39603024377847;
    "G54" = 97.97278799247276;
    "G55" = 97.98235259416693;
    "H56" = 97.51069335596023;
    "H57" = 97.4387859055598
}

foreach ($cellAddr in $updates.Keys) {
    $ws.Range($cellAddr).Value = $updates[$cellAddr]
}

$wb.Save()
